# Fix typo'd column headers (truncated "_ected" -> "_enacted")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "parental_leave_mandatory_enacted"
$ws.Range("C1").Value = "parental_leave_mandatory_not_yet_enacted"

# Fix truncated state names in column A (order matches original corrections)
$ws.Range("A28").Value = "Montana"
$ws.Range("A20").Value = "Louisiana"
$ws.Range("A35").Value = "North Carolina"
$ws.Range("A42").Value = "South Carolina"
$ws.Range("A16").Value = "Indiana"
